# "Error handling added on Excel"
#
# The "Settings" sheet stores a single hyperlinked "url" value (B2) that a
# UiPath process reads. We add a *second*, plain-text fallback URL (a
# differently-paginated variant of the same mobile.de search) in B2, and
# relocate the original hyperlinked URL out to G2 so both are available -
# the robot can fall back to the G2 link if the primary (B2) request fails.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember the original hyperlink's visible text (the long mobile.de URL
# currently sitting in B2) before we touch anything.
$originalUrl = $ws.Range("B2").Text

# Drop the hyperlink (and its relationship) that currently lives on B2.
$ws.Range("B2").Hyperlinks.Delete()

# B2 becomes a new, plain (non-hyperlinked) fallback URL value.
$ws.Range("B2").Value = "https://www.mobile.de/ro/automobil/mazda-cx-5/vhc:car,cnt:de,pgn:1,pgs:50,ms1:16800_33_,frn:2012,frx:2018,ful:diesel!electricity,mlx:100000"
$ws.Range("B2").Style = "Normal"

# Re-create the original hyperlink over on G2, pointing at the same
# external address, still showing the original URL text.
$ws.Hyperlinks.Add($ws.Range("G2"), $originalUrl)
$ws.Range("G2").Style = "Hyperlink"

# Match the author's final cursor position/selection.
$ws.Range("B11").Select()
